# "BOM For Me.xlsx" - BOM Report sheet:
#  - rows 11-13 (the line-item rows) no longer carry an explicit custom row
#    height; auto-fit them back to the sheet's default height.
#  - move the active selection to D15 (was C12).
#  - B8/C8 (=TODAY()/=NOW()) are volatile and simply get re-evaluated on
#    recalculation, which happens automatically on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BOM Report")

# Clear the explicit row height on rows 11:13 so Excel auto-fits them again
# (drops the ht="16.5" customHeight="1" attributes entirely).
$ws.Range("A11:A13").EntireRow.AutoFit()

# Move the selection to D15.
$ws.Activate()
$ws.Range("D15").Select()
